$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A5 must stay a literal text string "2025-08-20" (matching A2:A4 style),
# not get auto-converted into an Excel date serial number. Forcing the
# cell to Text format before assigning the value, then clearing the
# format afterwards, keeps the stored type as a string with no
# extra style applied (matching the existing date cells above it).
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "2025-08-20"
$ws.Range("A5").ClearFormats()

$ws.Range("B5").Value = 59.84999847412109
$ws.Range("C5").Value = 689.5999755859375
$ws.Range("D5").Value = 326.5499877929688
